$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 12.98179114390959
    "C2" = 11.04039597825339
    "D2" = 9.617039508785602
    "F2" = 34.34488192419382
    "G2" = 36.14186392745584
    "H2" = 16.28337400582205
    "J2" = 11.07483622179994
    "K2" = 8.967823268953866
    "L2" = 11.35246147243543
    "M2" = 15.28146341242424
    "N2" = 20.66716252841415
    "O2" = 25.73884412223041
    "B3" = 12.79266312756235
    "C3" = 11.05147512369043
    "D3" = 9.610161150389747
    "F3" = 34.42553427541914
    "G3" = 36.24062030855499
    "H3" = 16.32702328238774
    "J3" = 11.09548809078936
    "K3" = 8.819050637789925
    "L3" = 11.35914197096927
    "M3" = 15.25573976202018
    "N3" = 20.72181877517721
    "O3" = 25.8135559404136
    "B4" = 12.67720227071839
    "C4" = 11.05884087912727
    "D4" = 9.607245765369772
    "F4" = 34.48113570190502
    "G4" = 36.30922705170631
    "H4" = 16.35581503790987
    "J4" = 11.10889606137686
    "K4" = 8.727584534457685
    "L4" = 11.36442986316741
    "M4" = 15.24175831824633
    "N4" = 20.75699451078091
    "O4" = 25.86349063037778
    "B5" = 12.63037434821056
    "C5" = 11.06198448484832
    "D5" = 9.606388383132023
    "F5" = 34.5053220138001
    "G5" = 36.33918545745307
    "H5" = 16.36804907294456
    "J5" = 11.11454336018667
    "K5" = 8.690326724727319
    "L5" = 11.36688376579368
    "M5" = 15.23652093547063
    "N5" = 20.77173653728616
    "O5" = 25.88486060776303
    "B6" = 12.62261373822245
    "C6" = 11.0625150664414
    "D6" = 9.606266035350398
    "F6" = 34.50943041247418
    "G6" = 36.3442807511672
    "H6" = 16.37011080957014
    "J6" = 11.11549218317245
    "K6" = 8.684142357932467
    "L6" = 11.36730931863731
    "M6" = 15.23567918515969
    "N6" = 20.77420909250529
    "O6" = 25.88847075095545
    "B7" = 12.67656975449464
    "C7" = 11.05888269940289
    "D7" = 9.607232861346224
    "F7" = 34.48145570009006
    "G7" = 36.30962298572777
    "H7" = 16.35597800050096
    "J7" = 11.10897147944526
    "K7" = 8.727081940064117
    "L7" = 11.36446174550425
    "M7" = 15.24168581639747
    "N7" = 20.75719167493461
    "O7" = 25.86377469872667
    "B8" = 12.91647557424937
    "C8" = 11.04409943341611
    "D8" = 9.614397483666574
    "F8" = 34.3714284438545
    "G8" = 36.17425895037229
    "H8" = 16.29801137785948
    "J8" = 11.08180624845197
    "K8" = 8.916578170798186
    "L8" = 11.35451923965904
    "M8" = 15.27222025333965
    "N8" = 20.68567316230503
    "O8" = 25.76376179441028
    "B9" = 13.38974654987722
    "C9" = 11.01955988397999
    "D9" = 9.63874279919024
    "F9" = 34.2039397836633
    "G9" = 35.97219407228545
    "H9" = 16.20011244703339
    "J9" = 11.03428737196406
    "K9" = 9.285282767750171
    "L9" = 11.34439650468154
    "M9" = 15.34628623807266
    "N9" = 20.55820018406052
    "O9" = 25.59986418501351
    "B10" = 13.73593347779361
    "C10" = 11.0042198700053
    "D10" = 9.662787748115317
    "F10" = 34.11034653672864
    "G10" = 35.86256024095117
    "H10" = 16.13776949500539
    "J10" = 11.0028519480793
    "K10" = 9.551866078769333
    "L10" = 11.342624346399
    "M10" = 15.40907450819791
    "N10" = 20.47225928020675
    "O10" = 25.49909787081461
    "B11" = 13.89242672755144
    "C11" = 10.99782018210172
    "D11" = 9.67503559879596
    "F11" = 34.07416885200053
    "G11" = 35.82114734319421
    "H11" = 16.11148218882542
    "J11" = 10.98929975488158
    "K11" = 9.671697399622021
    "L11" = 11.34303672529786
    "M11" = 15.43939227788265
    "N11" = 20.43482163361909
    "O11" = 25.45752347306199
    "B12" = 13.95149105902379
    "C12" = 10.99547957819666
    "D12" = 9.679859236283498
    "F12" = 34.06138929834735
    "G12" = 35.80668392144935
    "H12" = 16.10182543435472
    "J12" = 10.98427496806794
    "K12" = 9.716827265363218
    "L12" = 11.34336703550958
    "M12" = 15.45111934637755
    "N12" = 20.42088212588336
    "O12" = 25.44239360628509
    "B13" = 13.93878006715594
    "C13" = 10.99597999172777
    "D13" = 9.678812167492845
    "F13" = 34.06410068103439
    "G13" = 35.80974464354011
    "H13" = 16.10389195609178
    "J13" = 10.98535238758312
    "K13" = 9.707119384183185
    "L13" = 11.34328816951551
    "M13" = 15.44858285205306
    "N13" = 20.42387370939979
    "O13" = 25.44562480746517
    "B14" = 13.89729022175265
    "C14" = 10.99762596155169
    "D14" = 9.675428737717574
    "F14" = 34.07309902957614
    "G14" = 35.8199329922158
    "H14" = 16.11068175814034
    "J14" = 10.98888421788777
    "K14" = 9.675415450795933
    "L14" = 11.34306041619913
    "M14" = 15.44035216803207
    "N14" = 20.43367007144948
    "O14" = 25.45626643425868
    "B15" = 13.87184936290158
    "C15" = 10.99864493989396
    "D15" = 9.673380383143611
    "F15" = 34.07873060720532
    "G15" = 35.82633242458713
    "H15" = 16.11487946286039
    "J15" = 10.99106150568774
    "K15" = 9.65596244745239
    "L15" = 11.34294355761172
    "M15" = 15.4353425415601
    "N15" = 20.43970150214954
    "O15" = 25.46286463316098
    "B16" = 13.72568171092949
    "C16" = 11.00464971450659
    "D16" = 9.662013449905839
    "F16" = 34.112839595998
    "G16" = 35.86543710649733
    "H16" = 16.13952910917701
    "J16" = 11.00375263038026
    "K16" = 9.544002340266076
    "L16" = 11.34262182605432
    "M16" = 15.40712798045289
    "N16" = 20.4747391817487
    "O16" = 25.50190068978728
    "B17" = 13.63572023532034
    "C17" = 11.00848136378186
    "D17" = 9.655373778758481
    "F17" = 34.13540324817342
    "G17" = 35.89159508121774
    "H17" = 16.15518148764261
    "J17" = 11.01172949052712
    "K17" = 9.474920332347883
    "L17" = 11.34273580454645
    "M17" = 15.39026458950118
    "N17" = 20.496657422033
    "O17" = 25.52694048963465
    "B18" = 13.58388681096914
    "C18" = 11.01073970653884
    "D18" = 9.651678259200976
    "F18" = 34.14898351616824
    "G18" = 35.90743665330486
    "H18" = 16.16437944805825
    "J18" = 11.01638799376321
    "K18" = 9.435053085945517
    "L18" = 11.34291612443835
    "M18" = 15.38073078213495
    "N18" = 20.50942026637602
    "O18" = 25.54174412474835
    "B19" = 13.56632313153918
    "C19" = 11.01151371100412
    "D19" = 9.650448300892482
    "F19" = 34.15368499318055
    "G19" = 35.91293701154896
    "H19" = 16.16752724450114
    "J19" = 11.01797739115896
    "K19" = 9.421533134148872
    "L19" = 11.34299692009142
    "M19" = 15.37753141409849
    "N19" = 20.5137683789948
    "O19" = 25.54682531631476
    "B20" = 13.64530649656483
    "C20" = 11.00806784249858
    "D20" = 9.656067826153421
    "F20" = 34.13293897683602
    "G20" = 35.88872809958558
    "H20" = 16.15349507323096
    "J20" = 11.01087305454393
    "K20" = 9.482288307012066
    "L20" = 11.3427118014735
    "M20" = 15.39204263312846
    "N20" = 20.4943080444014
    "O20" = 25.52423341480093
    "B21" = 13.90948254273349
    "C21" = 10.99714025563347
    "D21" = 9.676417516918818
    "F21" = 34.07043102433186
    "G21" = 35.81690733541044
    "H21" = 16.10867935289494
    "J21" = 10.9878439294213
    "K21" = 9.6847347010644
    "L21" = 11.34312259526509
    "M21" = 15.44276308145141
    "N21" = 20.43078620967307
    "O21" = 25.45312408012134
    "B22" = 14.08096922879327
    "C22" = 10.99048100417003
    "D22" = 9.6907979336219
    "F22" = 34.03494168706306
    "G22" = 35.77707281420899
    "H22" = 16.08112460448516
    "J22" = 10.97341733660683
    "K22" = 9.81558408335294
    "L22" = 11.34440562522727
    "M22" = 15.47734501966581
    "N22" = 20.39065382390213
    "O22" = 25.4102257723413
    "B23" = 13.98956773669434
    "C23" = 10.99399114228723
    "D23" = 9.683024873965817
    "F23" = 34.05339230975773
    "G23" = 35.79768256635489
    "H23" = 16.09567247556799
    "J23" = 10.98106009678338
    "K23" = 9.745893951242014
    "L23" = 11.34362838436593
    "M23" = 15.45875888918865
    "N23" = 20.41194703476923
    "O23" = 25.43279416927137
    "B24" = 13.64097289685179
    "C24" = 11.00825462264389
    "D24" = 9.655753668003353
    "F24" = 34.13405117967102
    "G24" = 35.89002176086933
    "H24" = 16.15425688066505
    "J24" = 11.01126002339133
    "K24" = 9.47895771312313
    "L24" = 11.34272229562191
    "M24" = 15.3912382770956
    "N24" = 20.49536969409144
    "O24" = 25.52545601270922
    "B25" = 13.26175138733153
    "C25" = 11.02572451079296
    "D25" = 9.631066130911449
    "F25" = 34.24407783131967
    "G25" = 36.02005207212563
    "H25" = 16.2249114244973
    "J25" = 11.04652979820712
    "K25" = 9.186115722159521
    "L25" = 11.34303672529786
    "M25" = 15.44858285205306
    "N25" = 20.42387370939979
    "O25" = 25.44562480746517
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

Write-Host "Updated $($values.Count) cells"